# Implements "Custom Excel Report" — builds two worksheets, Negative_Scenario
# and OnDemand_POC, each a small test-step report table, with a bold/double-
# underlined red "Script Failed at Step #N" caption at the bottom.

$wb = $excel.ActiveWorkbook

$headers = @("Test Step #", "Test Step Description", "Status")

# NOTE: the object handle returned by Worksheets.Add() tracks whatever sheet
# ends up at that slot rather than staying bound to the sheet it created
# (subsequent Add() calls can make it "follow" a different sheet), so we
# create + name both sheets first and then re-fetch stable handles via
# Worksheets.Item(name) before writing any cell data.
#
# Worksheets.Add() also inserts the new sheet immediately before the
# currently-active sheet, so "OnDemand_POC" is created first; once
# "Negative_Scenario" is added it becomes active and lands in front of
# "OnDemand_POC", producing the correct final tab order.

$wsNew2 = $wb.Worksheets.Add()
$wsNew2.Name = "OnDemand_POC"

$wsNew1 = $wb.Worksheets.Add()
$wsNew1.Name = "Negative_Scenario"

$ws1 = $wb.Worksheets.Item("Negative_Scenario")
$ws2 = $wb.Worksheets.Item("OnDemand_POC")

# ---------------------------------------------------------------------------
# Sheet 1: Negative_Scenario
# ---------------------------------------------------------------------------
$ws1.Range("A1").Value = $headers[0]
$ws1.Range("B1").Value = $headers[1]
$ws1.Range("C1").Value = $headers[2]

$ws1.Range("A2").Value = "Step 1"
$ws1.Range("B2").Value = " Entered User Name "
$ws1.Range("C2").Value = " Passed"

$ws1.Range("A3").Value = "Step 2"
$ws1.Range("B3").Value = " Entered Password "
$ws1.Range("C3").Value = " Passed"

$ws1.Range("A4").Value = ""

$ws1.Range("B5").Value = "Script Failed at Step #3"
$ws1.Range("B5").Font.Name = "Calibri"
$ws1.Range("B5").Font.Size = 15
$ws1.Range("B5").Font.Bold = $true
$ws1.Range("B5").Font.Underline = $true
$ws1.Range("B5").Font.Family = 4
$ws1.Range("B5").Font.Color = 723943

$ws1.Columns.Item(1).ColumnWidth = 10
$ws1.Columns.Item(2).ColumnWidth = 32
$ws1.Columns.Item(3).ColumnWidth = 15

# ---------------------------------------------------------------------------
# Sheet 2: OnDemand_POC
# ---------------------------------------------------------------------------
$ws2.Range("A1").Value = $headers[0]
$ws2.Range("B1").Value = $headers[1]
$ws2.Range("C1").Value = $headers[2]

# NOTE: nested array literals (an @() of @() pairs) get flattened by this
# host's PowerShell parser, so the step names/descriptions are kept as two
# parallel flat arrays and zipped by index instead.
$stepNames = @("Step 3", "Step 4", "Step 5", "Step 6", "Step 7", "Step 8", "Step 9")
$stepDescs = @(
    " Clicked SignIn Btn ",
    " Clicked Active Patient ",
    " Selected First Patient ",
    " Clicked On Demand under Action Menu ",
    " Clicked On Demand Campaign ",
    " Clicked Send Message Button ",
    " Validated Success Message "
)

$ws2.Range("A2").Value = "Step 1"
$ws2.Range("B2").Value = " Entered User Name "
$ws2.Range("C2").Value = " Passed"

$ws2.Range("A3").Value = "Step 2"
$ws2.Range("B3").Value = " Entered Password "
$ws2.Range("C3").Value = " Passed"

$row = 4
for ($i = 0; $i -lt $stepNames.Count; $i++) {
    $ws2.Cells.Item($row, 1).Value = $stepNames[$i]
    $ws2.Cells.Item($row, 2).Value = $stepDescs[$i]
    $ws2.Cells.Item($row, 3).Value = " Passed"
    $row = $row + 1
}

$ws2.Range("A11").Value = ""

$ws2.Range("B12").Value = "Script Failed at Step #10"
$ws2.Range("B12").Font.Name = "Calibri"
$ws2.Range("B12").Font.Size = 15
$ws2.Range("B12").Font.Bold = $true
$ws2.Range("B12").Font.Underline = $true
$ws2.Range("B12").Font.Family = 4
$ws2.Range("B12").Font.Color = 723943

$ws2.Columns.Item(1).ColumnWidth = 10
$ws2.Columns.Item(2).ColumnWidth = 32
$ws2.Columns.Item(3).ColumnWidth = 15

Write-Output "Negative_Scenario + OnDemand_POC sheets created"
